$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "{{ task }}" paragraph (currently paragraph 6): drop the
#    bookmarkStart/bookmarkEnd ("_GoBack") that used to sit at its end.
# ------------------------------------------------------------------
$pTask = $d.Paragraphs(6)
$taskXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="a4"/>' + `
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' + `
    '<w:tabs><w:tab w:val="left" w:pos="993"/></w:tabs>' + `
    '<w:rPr><w:lang w:val="en-US"/></w:rPr>' + `
  '</w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{{ task }}</w:t></w:r>' + `
'</w:p>'
$pTask.Range.InsertXML($taskXml)

# ------------------------------------------------------------------
# 2) "за период с ... по ..." paragraph (paragraph 3): re-split the
#    "за период с " run into word-by-word runs, wrap the date_from /
#    date_to placeholders with spellStart/spellEnd proofErr markers,
#    add lang="en-US" where the template now carries it, and move the
#    "_GoBack" bookmark here (to the very end of the paragraph).
# ------------------------------------------------------------------
$pPeriod = $d.Paragraphs(3)
$periodXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr>' + `
    '<w:spacing w:line="240" w:lineRule="auto"/>' + `
    '<w:ind w:firstLine="567"/>' + `
    '<w:jc w:val="center"/>' + `
    '<w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr>' + `
  '</w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>за</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>период</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>с</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>date</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>from</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> }} </w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>по</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> {{ </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>date</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>_</w:t></w:r>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>to</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
'</w:p>'
$pPeriod.Range.InsertXML($periodXml)

# ------------------------------------------------------------------
# 3) Remove the now-empty paragraph that used to follow it.
# ------------------------------------------------------------------
$pEmpty = $d.Paragraphs(4)
$pEmpty.Range.Delete()
